$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Per-row updates: Coin (B), Link (C), Price (D), Volume(1h) (E)
# Only the columns that actually changed are listed for each row.
$data = @{
    2  = @{ D = "43.876.48"; E = "  -0.30%  " }
    3  = @{ D = "2.300.94"; E = "  -1.30%  " }
    4  = @{ E = "  -0.18%  " }
    5  = @{ D = "101.41"; E = "  +3.96%  " }
    6  = @{ D = "270.59"; E = "  -0.50%  " }
    7  = @{ E = "  -0.14%  " }
    8  = @{ E = "  +0.00%  " }
    9  = @{ E = "  -3.00%  " }
    10 = @{ D = "45.00"; E = "  -2.74%  " }
    11 = @{ D = "0.0935"; E = "  -2.04%  " }
    12 = @{ E = "  -2.76%  " }
    13 = @{ E = "  +1.58%  " }
    14 = @{ D = "15.89"; E = "  +1.49%  " }
    15 = @{ D = "2.643.73"; E = "  -1.45%  " }
    16 = @{ E = "  -1.84%  " }
    17 = @{ D = "2.292.78"; E = "  -1.66%  " }
    18 = @{ D = "43.751.65"; E = "  -0.36%  " }
    19 = @{ E = "  +1.38%  " }
    20 = @{ D = "6.25"; E = "  -3.16%  " }
    22 = @{ E = "  +7.20%  " }
    23 = @{ D = "233.48"; E = "  -2.84%  " }
    24 = @{ D = "2.92"; E = "  +15.16%  " }
    25 = @{ E = "  -3.79%  " }
    26 = @{ D = "1.00"; E = "  +0.05%  " }
    27 = @{ E = "  -1.87%  " }
    28 = @{ D = "3.46"; E = "  -0.52%  " }
    29 = @{ D = "2.23"; E = "  -1.88%  " }
    30 = @{ B = "Monero"; C = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"; D = "177.30"; E = "  +2.03%  " }
    31 = @{ B = "InjectiveProtocol"; C = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"; D = "38.06"; E = "  -0.85%  " }
    32 = @{ D = "22.00"; E = "  -2.13%  " }
    33 = @{ D = "0.0894"; E = "  -1.47%  " }
    34 = @{ D = "5.47"; E = "  -0.90%  " }
    35 = @{ E = "  +0.53%  " }
    36 = @{ D = "4.76"; E = "  +6.85%  " }
    37 = @{ E = "  -1.36%  " }
    38 = @{ E = "  -2.87%  " }
    39 = @{ D = "3.54"; E = "  +4.37%  " }
    40 = @{ E = "  -1.31%  " }
    41 = @{ E = "  -1.26%  " }
    42 = @{ D = "1.39"; E = "  +0.99%  " }
    43 = @{ D = "12.21"; E = "  -0.73%  " }
    44 = @{ D = "64.79"; E = "  +3.64%  " }
    45 = @{ E = "  -3.99%  " }
    46 = @{ E = "  -3.31%  " }
    47 = @{ E = "  -1.38%  " }
    48 = @{ E = "  +0.65%  " }
    49 = @{ D = "98.70"; E = "  -1.77%  " }
    50 = @{ B = "Stacks"; C = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"; D = "1.54"; E = "  +11.53%  " }
    51 = @{ B = "WOONetwork"; C = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"; D = "0.441"; E = "  +5.48%  " }
}

foreach ($rowNum in $data.Keys) {
    $vals = $data[$rowNum]

    if ($vals.ContainsKey("B")) {
        $ws.Range("B$rowNum").Value = $vals["B"]
    }
    if ($vals.ContainsKey("C")) {
        $ws.Range("C$rowNum").Value = $vals["C"]
    }
    if ($vals.ContainsKey("D")) {
        # Prices are stored as plain text in this sheet (e.g. "43.876.48" or
        # "1.00"). Force text formatting so Excel doesn't reinterpret
        # numeric-looking values (like "101.41") as actual numbers, then
        # restore the default style so no stray formatting is left behind.
        $cell = $ws.Range("D$rowNum")
        $cell.NumberFormat = "@"
        $cell.Value = $vals["D"]
        $cell.Style = "Normal"
    }
    if ($vals.ContainsKey("E")) {
        $ws.Range("E$rowNum").Value = $vals["E"]
    }
}
